$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 11) mirroring row 10's pattern but one level deeper,
# filtering out non-positive numerator values.
$ws.Range("A11").Value = "Bottom"
$ws.Range("B11").Value = "Bottom"
$ws.Range("C11").Value = "Middle"
$ws.Range("D11").Value = "Middle"
$ws.Range("E11").Value = "Middle"
$ws.Range("F11").Value = "Eye"
$ws.Range("G11").Value = "Eye"
$ws.Range("H11").Value = "Eye"
$ws.Range("I11").Value = "Top"
$ws.Range("J11").Value = "Top"
$ws.Range("K11").Value = "Top"

$ws.Range("A11:K11").WrapText = $true

$null = $ws.Range("A11").Select()
